# Updates the acknowledgements on the last slide (slide 42) of the
# "class 5" Functions & PyGame deck: the final paragraph that used to
# credit Jeremy Pedersen's condensed 8-deck series is replaced with a
# new acknowledgement crediting the "Making Games with Python & Pygame"
# book by Al Sweigart (with a couple of italicised runs for the title).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item($p.Slides.Count)

# Locate the body textbox that holds the acknowledgements paragraphs
# (as opposed to the title placeholder or the logo picture).
$body = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
        if ($shp.TextFrame.TextRange.Text -like "*Slides condensed*") {
            $body = $shp
        }
    }
}

$tr = $body.TextFrame.TextRange

# Find the paragraph that needs replacing.
$target = $null
for ($i = 1; $i -le $tr.Paragraphs().Count; $i++) {
    $candidate = $tr.Paragraphs($i)
    if ($candidate.Text -like "*Slides condensed*") {
        $target = $candidate
    }
}

$newText = "This deck uses the style and formatting of Charles R. Severance’s slides, but some of the content and code is borrowed from the wonderful book Making Games with Python & Pygame by Al Sweigart. Like Charles Severance and Al Sweigart, I license these slides and all associated content under a CC license. "

$target.Text = $newText

# "Making Games with Python & " / "Pygame" are set in italics (two runs,
# since in the original deck "Pygame" is additionally flagged by the
# spell checker).
$target.Characters(144, 27).Font.Italic = $true
$target.Characters(171, 6).Font.Italic = $true

# Re-assert the (unchanged) font size on the remaining segments so they
# stay split into discrete runs, mirroring the original author's runs
# (" by Al ", "Sweigart", ". Like Charles Severance and Al ", "Sweigart",
# ", I license these slides and all associated content under a CC
# license. ").
$target.Characters(177, 7).Font.Size = 18
$target.Characters(184, 8).Font.Size = 18
$target.Characters(192, 32).Font.Size = 18
$target.Characters(224, 8).Font.Size = 18
$target.Characters(232, 72).Font.Size = 18
